# Generate Report for Handback
# Localization handback run: statuses flip from "Ready for handoff" to
# "Handed back: in sync with en-US", the handback timestamps are refreshed,
# and the stale "handback file is not latest" error details are cleared
# now that everything is back in sync.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: Status shown per-language for each source file ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("K2").Value = "2016-08-13 16:49:26"
$zhcn.Range("K3").Value = "2016-08-13 16:49:26"
$zhcn.Range("P2").Value = ""
$zhcn.Range("P3").Value = ""

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus
$dede.Range("K2").Value = "2016-08-13 16:49:37"
$dede.Range("K3").Value = "2016-08-13 16:49:37"
$dede.Range("P2").Value = ""
$dede.Range("P3").Value = ""

# --- Column widths: the longer status text / shorter error text reflow
# the "Status" and "Error Detail" columns on each sheet ---
$overview.Columns.Item(5).ColumnWidth = 29.166666667
$overview.Columns.Item(6).ColumnWidth = 29.166666667

$zhcn.Columns.Item(3).ColumnWidth = 29.166666667
$zhcn.Columns.Item(16).ColumnWidth = 12.833333333

$dede.Columns.Item(3).ColumnWidth = 29.166666667
$dede.Columns.Item(16).ColumnWidth = 12.833333333
